# Update column F ("dSF") values for rows 2-24 on Sheet1.
# Row 11 and row 25 are unchanged (remain 0) per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = -5
    3  = -2
    4  = 2
    5  = 4
    6  = 4
    7  = 2
    8  = -1
    9  = -3
    10 = 4
    12 = -2
    13 = -4
    14 = 2
    15 = -1
    16 = 3
    17 = -4
    18 = -5
    19 = -2
    20 = 3
    21 = 4
    22 = 3
    23 = -1
    24 = 2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
